$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 41440
$ws.Range("B13").Value = 2.5
$ws.Range("D13").Value = "Implementation tc12"

$ws.Range("A12").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 41442
$ws.Range("B14").Value = 2
$ws.Range("D14").Value = "Implementation tc12"

$ws.Range("A14").Select()
